$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row above row 2 (the existing data rows 2..41 shift
#    down to 3..42, carrying their values/styles with them).
# ------------------------------------------------------------------
$ws.Rows("2:2").Insert()

# ------------------------------------------------------------------
# 2. Populate the new row 2 with the latest circular entry.
# ------------------------------------------------------------------
$ws.Range("A2").Value = 41
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 278.75

# E2 holds a dd.mm.yyyy-looking label that must stay plain text (not get
# auto-parsed into a date serial) - force text entry, then drop back to
# the General format the other date-label cells use.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "09.10.2025"
$ws.Range("E2").NumberFormat = "General"

$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-09-october-2025.pdf"

# Match the number format / style already used by the other price cells.
$ws.Range("D2").NumberFormat = "0.000"

# ------------------------------------------------------------------
# 3. The row Insert() does not carry hyperlinks along with the shifted
#    rows, so rebuild the Circular Link hyperlinks from scratch for the
#    whole (now one-row-longer) list: F2:F22 each get a hyperlink whose
#    address matches the PDF link text already sitting in that cell.
#    (Rows 23 and below never had a hyperlink, before or after.)
# ------------------------------------------------------------------
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $ws.Hyperlinks.Item($i).Delete()
}

for ($r = 2; $r -le 22; $r++) {
    $target = $ws.Cells.Item($r, 6).Value
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
}
